$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 407, shifting rows 407:449 down to 408:450
$ws.Rows("407").Insert()

# Populate the new row 407 with the new weekly record
$ws.Range("A407").Value = 4
$ws.Range("B407").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C407").Value = "Los Lagos"
$ws.Range("D407").Value = 45194
$ws.Range("E407").Value = 10
$ws.Range("F407").Value = 100112044
$ws.Range("G407").Value = "Perejil"
$ws.Range("H407").Value = "Sin especificar"
$ws.Range("I407").Value = "Primera"
$ws.Range("J407").Value = 40
$ws.Range("K407").Value = 7000
$ws.Range("L407").Value = 7000
$ws.Range("M407").Value = 7000
$ws.Range("N407").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O407").Value = "Región Metropolitana"
$ws.Range("P407").Value = 2333
$ws.Range("Q407").Value = 3
$ws.Range("R407").Value = "Hortaliza"
